$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.688.51"
$ws.Range("E2").Value = "  -3.04%  "

# Row 3
$ws.Range("D3").Value = "2.094.66"
$ws.Range("E3").Value = "  -2.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "345.11"
$ws.Range("E5").Value = "  -2.16%  "

# Row 6
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").Value = "0.5145"
$ws.Range("E7").Value = "  -2.34%  "

# Row 8
$ws.Range("D8").Value = "0.4389"
$ws.Range("E8").Value = "  -3.79%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "52.54"
$ws.Range("E9").Value = "  -2.25%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.09264"
$ws.Range("E10").Value = "  +1.31%  "

# Row 11
$ws.Range("D11").Value = "1.171"
$ws.Range("E11").Value = "  -1.11%  "

# Row 12
$ws.Range("D12").Value = "24.84"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").Value = "2.096.56"
$ws.Range("E13").Value = "  -2.10%  "

# Row 14
$ws.Range("D14").Value = "8.277"
$ws.Range("E14").Value = "  +1.75%  "

# Row 15
$ws.Range("D15").Value = "6.748"
$ws.Range("E15").Value = "  -1.81%  "

# Row 16
$ws.Range("D16").Value = "99.53"
$ws.Range("E16").Value = "  -2.56%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001152"
$ws.Range("E17").Value = "  -1.73%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19
$ws.Range("D19").Value = "20.87"
$ws.Range("E19").Value = "  +7.02%  "

# Row 20
$ws.Range("D20").Value = "0.06648"
$ws.Range("E20").Value = "  -0.91%  "

# Row 21
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$ws.Range("D22").Value = "6.196"
$ws.Range("E22").Value = "  -2.35%  "

# Row 23
$ws.Range("D23").Value = "29.731.06"
$ws.Range("E23").Value = "  -3.14%  "

# Row 24
$ws.Range("D24").Value = "12.58"
$ws.Range("E24").Value = "  -2.04%  "

# Row 25
$ws.Range("D25").Value = "2.318"
$ws.Range("E25").Value = "  -2.68%  "

# Row 26
$ws.Range("D26").Value = "2.337.86"
$ws.Range("E26").Value = "  -1.42%  "

# Row 27
$ws.Range("D27").Value = "21.91"
$ws.Range("E27").Value = "  -2.60%  "

# Row 28
$ws.Range("E28").Value = "  -4.53%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.90"
$ws.Range("E29").Value = "  -1.74%  "

# Row 30
$ws.Range("D30").Value = "132.99"
$ws.Range("E30").Value = "  -2.51%  "

# Row 31
$ws.Range("D31").Value = "1.131"
$ws.Range("E31").Value = "  -7.53%  "

# Row 32
$ws.Range("E32").Value = "  -2.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.650"
$ws.Range("E33").Value = "  -1.03%  "

# Row 34
$ws.Range("D34").Value = "6.166"
$ws.Range("E34").Value = "  -3.15%  "

# Row 35
$ws.Range("D35").Value = "3.937"
$ws.Range("E35").Value = "  -1.93%  "

# Row 36
$ws.Range("D36").Value = "6.167"
$ws.Range("E36").Value = "  -0.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.30"
$ws.Range("E37").Value = "  -1.49%  "

# Row 38
$ws.Range("D38").Value = "0.02571"
$ws.Range("E38").Value = "  -3.14%  "

# Row 39
$ws.Range("D39").Value = "0.06704"
$ws.Range("E39").Value = "  -3.27%  "

# Row 40
$ws.Range("D40").Value = "12.45"
$ws.Range("E40").Value = "  -1.76%  "

# Row 41
$ws.Range("D41").Value = "0.6853"
$ws.Range("E41").Value = "  -1.89%  "

# Row 42
$ws.Range("D42").Value = "0.2221"
$ws.Range("E42").Value = "  -4.88%  "

# Row 43
$ws.Range("D43").Value = "1.298"
$ws.Range("E43").Value = "  +1.97%  "

# Row 44
$ws.Range("D44").Value = "0.6632"
$ws.Range("E44").Value = "  +2.62%  "

# Row 45
$ws.Range("D45").Value = "14.34"
$ws.Range("E45").Value = "  -2.88%  "

# Row 46
$ws.Range("D46").Value = "2.316"
$ws.Range("E46").Value = "  -1.46%  "

# Row 47
$ws.Range("E47").Value = "  -3.45%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000348"
$ws.Range("E48").Value = "  -6.13%  "

# Row 49
$ws.Range("D49").Value = "1.219"
$ws.Range("E49").Value = "  -2.87%  "

# Row 50
$ws.Range("D50").Value = "82.16"
$ws.Range("E50").Value = "  -1.04%  "

# Row 51
$ws.Range("D51").Value = "0.3304"
$ws.Range("E51").Value = "  +0.42%  "
